$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExactComparison")

# Update existing row labels (shared strings "Avergae Time"/"Avergae Iter" -> new
# "Average Time"/"Average Iter" strings)
$ws.Range("A2").Value = "Average Time"
$ws.Range("A3").Value = "Average Iter"

# Update the values for row 2 (Average Time)
$ws.Range("B2").Value = 0.021423909090909091
$ws.Range("C2").Value = 0.16751673636363631
$ws.Range("D2").Value = 1.0377829090909092

# Row 4 (Average fval) stays the same - no change needed

# Add new rows 5-8
$ws.Range("A5").Value = "Violation"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0

$ws.Range("A6").Value = "Average iter Bt"
$ws.Range("B6").Value = 0.097560975609756087
$ws.Range("C6").Value = 0.042553191489361701
$ws.Range("D6").Value = 0.018867924528301886

$ws.Range("A7").Value = "Average iter cg"
$ws.Range("B7").Value = 1.219512195121951
$ws.Range("C7").Value = 1.2127659574468088
$ws.Range("D7").Value = 1.1509433962264151

$ws.Range("A8").Value = "N converged"
$ws.Range("B8").Value = 11
$ws.Range("C8").Value = 11
$ws.Range("D8").Value = 11

# Adjust column widths to match target (engine quantizes ColumnWidth to 1/6
# character-unit steps, so these inputs are chosen to land as close as
# possible to the target stored widths of 13 / 14.5546875 / 13.5546875 / 13.5546875)
$ws.Columns.Item(1).ColumnWidth = 12.166666666666666
$ws.Columns.Item(2).ColumnWidth = 13.666666666666666
$ws.Columns.Item(3).ColumnWidth = 12.666666666666666
$ws.Columns.Item(4).ColumnWidth = 12.666666666666666
